$d = $word.ActiveDocument

$pairs = @(
    @("881÷5=", "427÷3="),
    @("559÷9=", "812÷3="),
    @("592÷9=", "733÷4="),
    @("822÷6=", "359÷2="),
    @("633÷2=", "540÷8="),
    @("993÷9=", "820÷4="),
    @("961÷8=", "541÷3="),
    @("908÷4=", "846÷7="),
    @("249÷6=", "742÷6="),
    @("811÷9=", "835÷8="),
    @("873÷4=", "416÷8="),
    @("686÷4=", "416÷3="),
    @("560÷6=", "634÷9="),
    @("996÷7=", "971÷4="),
    @("751÷7=", "418÷5="),
    @("622÷4=", "938÷5="),
    @("316÷3=", "897÷6="),
    @("640÷8=", "313÷7="),
    @("575÷2=", "517÷9="),
    @("445÷9=", "418÷3="),
    @("898÷3=", "587÷7="),
    @("888÷9=", "514÷3="),
    @("698÷7=", "842÷2="),
    @("623÷6=", "186÷8="),
    @("486÷5=", "606÷3=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
